$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows for species added to the fit-item template (PEs: monophosphorylated
# PER, total PER, PER mRNA, nuclear PER, EmptySet, unphosphorylated PER,
# biphosphorylated PER), appended below the existing reaction-parameter rows.
$rows = @(
    @("monophosphorylated PER", "0.25", "1e-06", "1000000", "reactions",  "Species", "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[CYTOPLASM],Vector=Metabolites[monophosphorylated PER]"),
    @("total PER",              "1.0",  "1e-06", "1000000", "assignment", "Species", "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[CYTOPLASM],Vector=Metabolites[total PER]"),
    @("PER mRNA",               "0.1",  "1e-06", "1000000", "reactions",  "Species", "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[CYTOPLASM],Vector=Metabolites[PER mRNA]"),
    @("nuclear PER",            "0.25", "1e-06", "1000000", "reactions",  "Species", "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[NUCLEUS],Vector=Metabolites[nuclear PER]"),
    @("EmptySet",               "0.0",  "1e-06", "1000000", "fixed",      "Species", "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[default],Vector=Metabolites[EmptySet]"),
    @("unphosphorylated PER",   "0.25", "1e-06", "1000000", "reactions",  "Species", "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[CYTOPLASM],Vector=Metabolites[unphosphorylated PER]"),
    @("biphosphorylated PER",   "0.25", "1e-06", "1000000", "reactions",  "Species", "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[CYTOPLASM],Vector=Metabolites[biphosphorylated PER]")
)

$firstRow = 20
$lastRow = 26

# Column A picks up the same bold/bordered/centered style already used by the
# existing parameter-name column - copy the format down from the row above.
$ws.Range("A" + ($firstRow - 1)).Copy() | Out-Null
$ws.Range("A" + $firstRow + ":A" + $lastRow).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Force textual (shared-string) storage on B:G for the whole new block so the
# numeric-looking values (0.25, 1e-06, 1000000, ...) aren't auto-converted to
# real numbers by Excel, then strip the helper number format back off again
# so no extra cell styling lingers on the written cells.
$dataRange = $ws.Range("B" + $firstRow + ":G" + $lastRow)
$dataRange.NumberFormat = "@"

$r = $firstRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}

$dataRange.ClearFormats()
